$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The upstream odds feed re-ordered several pairs of fixtures that share the
# same kickoff date/time. Swap every column except "A" (the sheet's own
# sequential row index) between each pair of rows so the row position stays
# put while the match data trades places.
$swapPairs = @(
    @(12, 13),
    @(24, 25),
    @(35, 36),
    @(37, 38),
    @(46, 47),
    @(70, 72),
    @(132, 133)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")
    $tmp = $range1.Value2
    $range1.Value2 = $range2.Value2
    $range2.Value2 = $tmp
}

# Remove the six not-yet-played fixtures (rows 142-147) that dropped out of
# the feed entirely.
$ws.Range("A142:AC147").EntireRow.Delete()
